{"js": "// Add two new paragraphs (\"Lab 3\" and \"Block 2\") at the end of the document\n// body - right after the existing \"Software Configuration Management\"\n// paragraph and before the section break - matching the target diff.\nconst body = context.document.body;\n\nbody.insertParagraph(\"Lab 3\", Word.InsertLocation.end);\nbody.insertParagraph(\"Block 2\", Word.InsertLocation.end);\n\nawait context.sync();\n", "ps1": "# Add two new paragraphs (\"Lab 3\" and \"Block 2\") at the end of the document,\n# right after the existing \"Software Configuration Management\" paragraph and\n# before the section break - matching the target diff.\n$d = $word.ActiveDocument\n\n$p1 = $d.Paragraphs.Add()\n$p1.Range.Text = \"Lab 3\"\n\n$p2 = $d.Paragraphs.Add()\n$p2.Range.Text = \"Block 2\"\n"}
